$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 44
$ws1.Range("F3").Value = 152
$ws1.Range("F4").Value = 65
$ws1.Range("F6").Value = 1515
$ws1.Range("F7").Value = 1037
$ws1.Range("F9").Value = 212
$ws1.Range("F10").Value = 151
$ws1.Range("F13").Value = 183

# Sheet "全部类型" (sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 44
$ws4.Range("F3").Value = 152
$ws4.Range("F4").Value = 65
$ws4.Range("F6").Value = 1515
$ws4.Range("F8").Value = 1037
$ws4.Range("F10").Value = 212
$ws4.Range("F11").Value = 151
$ws4.Range("F14").Value = 183
